$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d7 = @"
<b>ASPECTOS HISTÓRICOS E EVOLUTIVOS</b>
do:
<ul>
	<li><b>MRP I</b></li>
	<li><b>MRP II</b></li>
	<li><b>ERP I</b></li>
	<li><b>ERP II</b></li>
	<li><b>ERPs baseados em nuvem (SaaS)</b></li>
</ul>

"@

$e7 = @"
<ul>
	<li><b>MRP I:</b>
<ul>
	<li>Foco no planejamento de materiais necessários para a produção;</li>
	<li>Estrutura baseada em listas de materiais (BOM – Bill of Materials), programação de ordens de</li>
	<li>compra e produção;</li>
	<li>Objetivo central: minimizar estoques e garantir disponibilidade de insumos.</li>
</ul>
</li>
	<li><b>MRP II</b>
<ul>
	<li>Expansão do MRP para incluir capacidade produtiva, finanças e simulações;</li>
	<li>Incorporação de planejamento financeiro integrado com o plano mestre de produção;</li>
	<li>Primeiro movimento em direção à integração interfuncional.</li>
</ul>
</li>
	<li><b>ERP I</b>
<ul>
	<li>Integração total de todos os departamentos e funções empresariais;</li>
	<li>Inclusão de módulos de vendas, distribuição, contabilidade, RH, manutenção e outros;</li>
	<li>Capacidade de operar em ambientes multiempresa e multinacionais.</li>
</ul>
</li>
	<li><b>ERP II</b>
<ul>
	<li>com integração entre empresas via internet e cadeia de suprimentos (SCM).</li>
</ul>
</li>
	<li><b>ERPs baseados em nuvem (SaaS)</b>
<ul>
	<li>que ampliam escalabilidade e flexibilidade.</li>
</ul>
</li>
</ul>
"@

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(7, 3).Value = "Gestão da Produção e Operações"
$ws.Cells.Item(7, 4).Value = $d7
$ws.Cells.Item(7, 5).Value = $e7
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
